$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '64.241.91'
$ws.Range('D2').Style = "Normal"

$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('E2').Style = "Normal"

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.491.97'
$ws.Range('D3').Style = "Normal"

$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('E3').Style = "Normal"

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E4').Style = "Normal"

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '587.38'
$ws.Range('D5').Style = "Normal"

$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.21%  '
$ws.Range('E5').Style = "Normal"

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '134.15'
$ws.Range('D6').Style = "Normal"

$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.67%  '
$ws.Range('E6').Style = "Normal"

$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E8').Style = "Normal"

$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -0.35%  '
$ws.Range('E9').Style = "Normal"

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.25'
$ws.Range('D10').Style = "Normal"

$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +1.93%  '
$ws.Range('E10').Style = "Normal"

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.386'
$ws.Range('D11').Style = "Normal"

$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +1.81%  '
$ws.Range('E11').Style = "Normal"

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.085.84'
$ws.Range('D12').Style = "Normal"

$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -1.01%  '
$ws.Range('E12').Style = "Normal"

$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +1.04%  '
$ws.Range('E13').Style = "Normal"

$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +1.33%  '
$ws.Range('E14').Style = "Normal"

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.506.55'
$ws.Range('D15').Style = "Normal"

$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -0.75%  '
$ws.Range('E15').Style = "Normal"

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '25.78'
$ws.Range('D16').Style = "Normal"

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '64.329.24'
$ws.Range('D17').Style = "Normal"

$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('E17').Style = "Normal"

$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range('E18').Style = "Normal"

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '5.74'
$ws.Range('D19').Style = "Normal"

$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +2.24%  '
$ws.Range('E19').Style = "Normal"

$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -3.07%  '
$ws.Range('E20').Style = "Normal"

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '394.10'
$ws.Range('D21').Style = "Normal"

$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +2.34%  '
$ws.Range('E21').Style = "Normal"

$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.81%  '
$ws.Range('E22').Style = "Normal"

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.631.75'
$ws.Range('D23').Style = "Normal"

$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.97%  '
$ws.Range('E23').Style = "Normal"

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '74.73'
$ws.Range('D24').Style = "Normal"

$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +0.98%  '
$ws.Range('E24').Style = "Normal"

$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('E25').Style = "Normal"

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '5.73'
$ws.Range('D26').Style = "Normal"

$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('E26').Style = "Normal"

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000116'
$ws.Range('D27').Style = "Normal"

$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +0.10%  '
$ws.Range('E27').Style = "Normal"

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.993'
$ws.Range('D28').Style = "Normal"

$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -0.75%  '
$ws.Range('E28').Style = "Normal"

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.34'
$ws.Range('D29').Style = "Normal"

$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -2.02%  '
$ws.Range('E29').Style = "Normal"

$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('E30').Style = "Normal"

$ws.Range('B31').NumberFormat = "@"
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('B31').Style = "Normal"

$ws.Range('C31').NumberFormat = "@"
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C31').Style = "Normal"

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.21'
$ws.Range('D31').Style = "Normal"

$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -2.33%  '
$ws.Range('E31').Style = "Normal"

$ws.Range('B32').NumberFormat = "@"
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('B32').Style = "Normal"

$ws.Range('C32').NumberFormat = "@"
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('C32').Style = "Normal"

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.48'
$ws.Range('D32').Style = "Normal"

$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -5.47%  '
$ws.Range('E32').Style = "Normal"

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.513.87'
$ws.Range('D33').Style = "Normal"

$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.59%  '
$ws.Range('E33').Style = "Normal"

$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +3.54%  '
$ws.Range('E34').Style = "Normal"

$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E35').Style = "Normal"

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '23.40'
$ws.Range('D36').Style = "Normal"

$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.77%  '
$ws.Range('E36').Style = "Normal"

$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -4.40%  '
$ws.Range('E37').Style = "Normal"

$ws.Range('B38').NumberFormat = "@"
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('B38').Style = "Normal"

$ws.Range('C38').NumberFormat = "@"
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('C38').Style = "Normal"

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.55'
$ws.Range('D38').Style = "Normal"

$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('E38').Style = "Normal"

$ws.Range('B39').NumberFormat = "@"
$ws.Range('B39').Value = 'Aptos'
$ws.Range('B39').Style = "Normal"

$ws.Range('C39').NumberFormat = "@"
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('C39').Style = "Normal"

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.88'
$ws.Range('D39').Style = "Normal"

$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -0.65%  '
$ws.Range('E39').Style = "Normal"

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '166.64'
$ws.Range('D40').Style = "Normal"

$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +3.49%  '
$ws.Range('E40').Style = "Normal"

$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -1.12%  '
$ws.Range('E42').Style = "Normal"

$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('E43').Style = "Normal"

$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -4.46%  '
$ws.Range('E44').Style = "Normal"

$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.78%  '
$ws.Range('E45').Style = "Normal"

$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +2.01%  '
$ws.Range('E46').Style = "Normal"

$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -3.88%  '
$ws.Range('E47').Style = "Normal"

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.460.97'
$ws.Range('D48').Style = "Normal"

$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -0.45%  '
$ws.Range('E48').Style = "Normal"

$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -0.82%  '
$ws.Range('E49').Style = "Normal"

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.893'
$ws.Range('D50').Style = "Normal"

$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -1.84%  '
$ws.Range('E50').Style = "Normal"

$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -1.26%  '
$ws.Range('E51').Style = "Normal"
